$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> new value for column F (dSF)
$changes = @{
    3  = -2
    13 = -1
    14 = -2
    15 = 0
    17 = -3
    21 = 1
    22 = -4
    33 = 5
    40 = 6
    42 = -2
    47 = -2
    49 = -4
    54 = -2
    55 = -4
    56 = -7
    57 = -6
    58 = 1
    59 = 6
    60 = -1
    64 = -2
}

foreach ($row in $changes.Keys) {
    $ws.Range("F$row").Value = $changes[$row]
}
